$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 5.682516
$ws.Range("H2").Value = 17.047548
$ws.Range("I2").Value = 0.4522589164991918
$ws.Range("J2").Value = 0.4522589164991919
$ws.Range("M2").Value = 0.2901893333333334
$ws.Range("N2").Value = 0.870568
$ws.Range("O2").Value = 0.03429389578125064
$ws.Range("P2").Value = 0.03429389578125064
$ws.Range("Q2").Value = 1.649005529696
$ws.Range("R2").Value = 14.841049767264
$ws.Range("S2").Value = 0.01550972014856462
$ws.Range("T2").Value = 0.01550972014856462

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 5.682516
$ws.Range("H3").Value = 17.047548
$ws.Range("I3").Value = 0.4522589164991918
$ws.Range("J3").Value = 0.4522589164991919
$ws.Range("O3").Value = 0.8402845891331153
$ws.Range("P3").Value = 0.8402845891331153
$ws.Range("Q3").Value = 40.404681428944
$ws.Range("R3").Value = 363.642132860496
$ws.Range("S3").Value = 0.3800261978323113
$ws.Range("T3").Value = 0.3800261978323113

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 5.682516
$ws.Range("H4").Value = 17.047548
$ws.Range("I4").Value = 0.4522589164991918
$ws.Range("J4").Value = 0.4522589164991919
$ws.Range("O4").Value = 0.1254215150856341
$ws.Range("P4").Value = 0.1254215150856341
$ws.Range("Q4").Value = 6.030833394908
$ws.Range("R4").Value = 54.277500554172
$ws.Range("S4").Value = 0.05672299851831593
$ws.Range("T4").Value = 0.05672299851831592

# Row 5
$ws.Range("I5").Value = 0.4336933920535619
$ws.Range("J5").Value = 0.433693392053562
$ws.Range("M5").Value = 0.2901893333333334
$ws.Range("N5").Value = 0.870568
$ws.Range("O5").Value = 0.03429389578125064
$ws.Range("P5").Value = 0.03429389578125064
$ws.Range("Q5").Value = 1.58131277372
$ws.Range("R5").Value = 14.23181496348
$ws.Range("S5").Value = 0.01487303598810193
$ws.Range("T5").Value = 0.01487303598810193

# Row 6
$ws.Range("I6").Value = 0.4336933920535619
$ws.Range("J6").Value = 0.433693392053562
$ws.Range("O6").Value = 0.8402845891331153
$ws.Range("P6").Value = 0.8402845891331153
$ws.Range("S6").Value = 0.3644258737514744
$ws.Range("T6").Value = 0.3644258737514744

# Row 7
$ws.Range("I7").Value = 0.4336933920535619
$ws.Range("J7").Value = 0.433693392053562
$ws.Range("O7").Value = 0.1254215150856341
$ws.Range("P7").Value = 0.1254215150856341
$ws.Range("Q7").Value = 5.783263737935001
$ws.Range("S7").Value = 0.05439448231398565
$ws.Range("T7").Value = 0.05439448231398565

# Row 8
$ws.Range("I8").Value = 0.1140476914472462
$ws.Range("J8").Value = 0.1140476914472462
$ws.Range("M8").Value = 0.2901893333333334
$ws.Range("N8").Value = 0.870568
$ws.Range("O8").Value = 0.03429389578125064
$ws.Range("P8").Value = 0.03429389578125064
$ws.Range("Q8").Value = 0.4158354141502222
$ws.Range("R8").Value = 3.742518727352
$ws.Range("S8").Value = 0.003911139644584091
$ws.Range("T8").Value = 0.003911139644584092

# Row 9
$ws.Range("I9").Value = 0.1140476914472462
$ws.Range("J9").Value = 0.1140476914472462
$ws.Range("O9").Value = 0.8402845891331153
$ws.Range("P9").Value = 0.8402845891331153
$ws.Range("S9").Value = 0.09583251754932957
$ws.Range("T9").Value = 0.0958325175493296

# Row 10
$ws.Range("I10").Value = 0.1140476914472462
$ws.Range("J10").Value = 0.1140476914472462
$ws.Range("O10").Value = 0.1254215150856341
$ws.Range("P10").Value = 0.1254215150856341
$ws.Range("Q10").Value = 1.520816065974556
$ws.Range("S10").Value = 0.01430403425333253
$ws.Range("T10").Value = 0.01430403425333253
